$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "Transportation"
$ws.Range("C18").Value = "Cable"
$ws.Range("D18").Value = "Single role"
$ws.Range("E18").Value = "Distributed"
$ws.Range("F18").Value = "Agent positions in formation and common velocity set at beginning"
$ws.Range("G18").Value = "Passivity-based PD (internal feedback) +  feedback control (formation)"
$ws.Range("H18").Value = "Agents' state"

$ws.Range("A18:F18").Style = $ws.Range("A17:F17").Style
$ws.Range("G18:H18").Style = $ws.Range("G17:H17").Style

$ws.Range("A18").Select()
